# Insert a new daily price record as row 9 (Vega Monumental Concepción - Ajo,
# China/Primera, fecha 2022-01-08 i.e. serial 44552), shifting every
# subsequent row down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()

$ws.Cells.Item(9, 1).Value2 = 11
$ws.Cells.Item(9, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value2 = "Bíobío"
$ws.Cells.Item(9, 4).Value2 = 44552
$ws.Cells.Item(9, 5).Value2 = 8
$ws.Cells.Item(9, 6).Value2 = 100112003
$ws.Cells.Item(9, 7).Value2 = "Ajo"
$ws.Cells.Item(9, 8).Value2 = "Chino"
$ws.Cells.Item(9, 9).Value2 = "Primera"
$ws.Cells.Item(9, 10).Value2 = 400
$ws.Cells.Item(9, 11).Value2 = 17000
$ws.Cells.Item(9, 12).Value2 = 18000
$ws.Cells.Item(9, 13).Value2 = 17500
$ws.Cells.Item(9, 14).Value2 = "`$/caja 10 kilos"
$ws.Cells.Item(9, 15).Value2 = "China"
$ws.Cells.Item(9, 16).Value2 = 1750
$ws.Cells.Item(9, 17).Value2 = 10
$ws.Cells.Item(9, 18).Value2 = "Hortaliza"
